$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 0.3326388888888889
$ws.Range("F10").Value = 0.4993055555555556
$ws.Range("E11").Value = 0.5847222222222223
$ws.Range("F11").Value = 0.7597222222222222
$ws.Range("E12").Value = 0.3354166666666666
$ws.Range("F12").Value = 0.5
$ws.Range("E13").Value = 0.5854166666666667
$ws.Range("F13").Value = 0.7472222222222222
$ws.Range("E14").Value = 0.3305555555555555
$ws.Range("F14").Value = 0.5020833333333333
$ws.Range("E15").Value = 0.5798611111111112
$ws.Range("F15").Value = 0.7506944444444444
$ws.Range("E16").Value = 0.3340277777777778
$ws.Range("F16").Value = 0.4979166666666667
$ws.Range("E17").Value = 0.5888888888888889
$ws.Range("E18").Value = 0.3347222222222222
$ws.Range("F18").Value = 0.5027777777777778
$ws.Range("E19").Value = 0.58125
$ws.Range("F19").Value = 0.7513888888888889
$ws.Range("E20").Value = 0.33125
$ws.Range("F20").Value = 0.4993055555555556
$ws.Range("E21").Value = 0.58125
$ws.Range("F21").Value = 0.7569444444444444
$ws.Range("E22").Value = 0.3347222222222222
$ws.Range("F22").Value = 0.5006944444444444
$ws.Range("E23").Value = 0.5930555555555556
$ws.Range("F23").Value = 0.7506944444444444
$ws.Range("E24").Value = 0.33125
$ws.Range("F24").Value = 0.4972222222222222
$ws.Range("E25").Value = 0.5895833333333333
$ws.Range("F25").Value = 0.7493055555555556
$ws.Range("E26").Value = 0.33125
$ws.Range("F26").Value = 0.5
$ws.Range("E27").Value = 0.58125
$ws.Range("F27").Value = 0.7569444444444444
$ws.Range("E28").Value = 0.3333333333333333
$ws.Range("F28").Value = 0.5013888888888889
$ws.Range("E29").Value = 0.5798611111111112
$ws.Range("F29").Value = 0.7493055555555556
$ws.Range("E30").Value = 0.3347222222222222
$ws.Range("E31").Value = 0.5798611111111112
$ws.Range("F31").Value = 0.7486111111111111
$ws.Range("F32").Value = 0.5027777777777778
$ws.Range("E33").Value = 0.5868055555555556
$ws.Range("F33").Value = 0.7479166666666667
$ws.Range("E34").Value = 0.3305555555555555
$ws.Range("F34").Value = 0.4986111111111111
$ws.Range("E35").Value = 0.5819444444444445
$ws.Range("F35").Value = 0.7520833333333333
$ws.Range("E36").Value = 0.3319444444444444
$ws.Range("F36").Value = 0.5027777777777778
$ws.Range("E37").Value = 0.5881944444444445
$ws.Range("F37").Value = 0.7555555555555555
$ws.Range("F38").Value = 0.5006944444444444
$ws.Range("F39").Value = 0.7472222222222222
$ws.Range("E40").Value = 0.3368055555555556
$ws.Range("F40").Value = 0.5034722222222222
$ws.Range("F41").Value = 0.7479166666666667
$ws.Range("E42").Value = 0.3340277777777778
$ws.Range("F42").Value = 0.4993055555555556
$ws.Range("E43").Value = 0.5909722222222222
$ws.Range("F43").Value = 0.7576388888888889
$ws.Range("E44").Value = 0.3319444444444444
$ws.Range("F44").Value = 0.4993055555555556
$ws.Range("E45").Value = 0.5881944444444445
$ws.Range("F45").Value = 0.7520833333333333
$ws.Range("E46").Value = 0.3333333333333333
$ws.Range("F46").Value = 0.5027777777777778
$ws.Range("E47").Value = 0.5798611111111112
$ws.Range("F47").Value = 0.7534722222222222
$ws.Range("E48").Value = 0.33125
$ws.Range("F48").Value = 0.4972222222222222
$ws.Range("E49").Value = 0.5868055555555556
$ws.Range("F49").Value = 0.7479166666666667
$ws.Range("E50").Value = 0.3354166666666666
$ws.Range("F50").Value = 0.5020833333333333
$ws.Range("E51").Value = 0.5805555555555556
$ws.Range("F51").Value = 0.7486111111111111
$ws.Range("F52").Value = 0.5
$ws.Range("E53").Value = 0.5868055555555556
$ws.Range("F53").Value = 0.75
$ws.Range("E54").Value = 0.3319444444444444
$ws.Range("F54").Value = 0.5006944444444444
$ws.Range("E55").Value = 0.5819444444444445
$ws.Range("F55").Value = 0.75
$ws.Range("E56").Value = 0.3361111111111111
$ws.Range("F56").Value = 0.4965277777777778
$ws.Range("E57").Value = 0.5916666666666667
$ws.Range("F57").Value = 0.7472222222222222
$ws.Range("E58").Value = 0.3333333333333333
$ws.Range("F58").Value = 0.5020833333333333
$ws.Range("E59").Value = 0.5819444444444445
$ws.Range("F59").Value = 0.7506944444444444
$ws.Range("E60").Value = 0.3333333333333333
$ws.Range("F60").Value = 0.4972222222222222
$ws.Range("E61").Value = 0.5819444444444445
$ws.Range("F61").Value = 0.7576388888888889
$ws.Range("E62").Value = 0.33125
$ws.Range("F62").Value = 0.4965277777777778
